$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "MIE(8.337755133772903, 5.0572030805735935, -10.11305051265931, 11.922641962745328)"
$ws.Range("C2").Value = "NIG(2.1910975699432873, 1.7816554348426665, 3.0847467423744863, 6.889555637008751)"
$ws.Range("D2").Value = "NIG(0.7930632690498539, 0.5094026695405718, 1.419497209597988, 3.0270772577669702)"
$ws.Range("E2").Value = "NIG(1.0805655219587362, 0.8308985024863275, 3.537203511217257, 5.826580913183122)"
